$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.395667552947998
$ws.Range("B1").Value = 3.036719083786011
$ws.Range("C1").Value = 2.712340593338013
$ws.Range("D1").Value = 1.530779242515564
$ws.Range("E1").Value = 1.160770177841187
